# Many hard things added:
#  - Sheet1 gets a new "Id" column inserted before the existing "Name" column
#  - A new "Sheet2" is added (after Sheet1) with "Id"/"Age" columns

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new sheet right after Sheet1 first (while Sheet1's layout is
# still untouched) so we control the shared-string insertion order:
# "Age" needs to land in the shared-strings table before "Id".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Sheet2: Id / Age ---
$ws2.Range("B1").Value = "Age"
$ws2.Range("A1").Value = "Id"
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 34
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 35
$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 15
$ws2.Range("E14").Select()

# --- Sheet1: insert a new Id column before the existing Name column ---
$ws1.Columns.Item(1).Insert()
$ws1.Range("A1").Value = "Id"
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 2
$ws1.Range("A4").Value = 3
$ws1.Range("A5").Value = 4
$ws1.Range("D8").Select()
